$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 75

# Column A needs to hold the text "01-04-2021" as a shared string (like the
# other period labels in column A), not have Excel auto-convert it into a
# date serial number. Using a TEXT() formula in a scratch cell and pasting
# the result back as a value forces a plain text cell without leaving any
# extra/unused cell styles behind (unlike toggling NumberFormat directly).
$scratch = $ws.Range("Z1")
$scratch.Formula = '=TEXT("01-04-2021","@")'
$scratch.Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("B" + $newRow).Value = 41591.51
$ws.Range("C" + $newRow).Value = 11810.18
$ws.Range("D" + $newRow).Value = 29781.34
